$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates. D-column price strings are forced to text via a
# leading apostrophe so numeric-looking values (e.g. '1.00', '0.919')
# are stored as text, matching the source inlineStr cell type instead of
# being auto-converted to numbers by Excel's type inference.

$ws.Range("D2").Value = "'66.496.99"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "'3.345.92"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'190.89"
$ws.Range("E5").Value = "  +4.94%  "
$ws.Range("D6").Value = "'565.10"
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.590"
$ws.Range("E8").Value = "  -1.23%  "
$ws.Range("D9").Value = "'3.333.87"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").Value = "'0.186"
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("D11").Value = "'0.591"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").Value = "'47.99"
$ws.Range("E12").Value = "  +0.91%  "
$ws.Range("D13").Value = "'0.0000272"
$ws.Range("E13").Value = "  +2.34%  "
$ws.Range("D14").Value = "'8.72"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "'3.879.82"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").Value = "'612.50"
$ws.Range("E16").Value = "  +1.60%  "
$ws.Range("D17").Value = "'18.20"
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("D18").Value = "'66.589.15"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").Value = "'3.347.92"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").Value = "'11.20"
$ws.Range("E21").Value = "  -2.06%  "
$ws.Range("D22").Value = "'0.919"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").Value = "'18.49"
$ws.Range("E23").Value = "  +9.88%  "
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("D25").Value = "'101.63"
$ws.Range("E25").Value = "  +2.28%  "
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").Value = "'2.76"
$ws.Range("E27").Value = "  +2.44%  "
$ws.Range("D28").Value = "'9.78"
$ws.Range("E28").Value = "  +4.99%  "
$ws.Range("D29").Value = "'8.72"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "'30.56"
$ws.Range("E30").Value = "  +0.64%  "
$ws.Range("E31").Value = "  +9.14%  "
$ws.Range("D32").Value = "'4.11"
$ws.Range("E32").Value = "  +9.14%  "
$ws.Range("D33").Value = "'572.10"
$ws.Range("E33").Value = "  +2.02%  "
$ws.Range("D34").Value = "'11.17"
$ws.Range("E34").Value = "  +1.02%  "
$ws.Range("E35").Value = "  +1.65%  "
$ws.Range("D36").Value = "'3.752.96"
$ws.Range("E36").Value = "  -2.13%  "
$ws.Range("D37").Value = "'57.45"
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").Value = "'0.0₃0733"
$ws.Range("E39").Value = "  +2.82%  "
$ws.Range("D40").Value = "'34.33"
$ws.Range("E40").Value = "  +7.03%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.132"
$ws.Range("E41").Value = "  +3.91%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'3.33"
$ws.Range("E42").Value = "  -2.08%  "
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "'2.74"
$ws.Range("E43").Value = "  +3.89%  "
$ws.Range("B44").Value = "CoreDAO"
$ws.Range("C44").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D44").Value = "'3.41"
$ws.Range("E44").Value = "  +1.89%  "
$ws.Range("D45").Value = "'0.344"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").Value = "'0.0426"
$ws.Range("E46").Value = "  +3.01%  "
$ws.Range("D47").Value = "'3.24"
$ws.Range("E47").Value = "  +3.14%  "
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("E51").Value = "  +2.80%  "
